# Update the "Test Data" sheet: clear the (now stale) query-table row 23
# and move the selection/scroll position, mirroring a manual edit made in
# Excel after the underlying data connection dropped that row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Data")
$ws.Select()

# Row 23 previously held a full record (BOM/Test/Vendor/dates/result/status).
# Clear all of its values while leaving the date-formatted cells' number
# formatting intact, so it becomes an empty row just like row 24 below it.
$ws.Range("A23:H23").ClearContents()

# Reflect the user's new selection/scroll position on the sheet.
$ws.Range("A23:Z23").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 13
$win.ScrollColumn = 7
